$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker row for the "TM" expert response collection (Oct-24-2023).
# Column layout: A=Task_type, B=Date, C=expert, D=submission_file_name, E=response_collected
$ws.Range("A16").Value = "pairwise"
$ws.Range("C16").Value = "TM"

# "Oct-24-2023" looks like a real date, so a plain Range.Value assignment would
# get auto-converted into a date serial number (and would also force a brand
# new number-format style to be created). To keep it as plain text - matching
# how the original workbook stores it as a shared string with no cell style -
# build it as a text formula on a scratch cell, then copy/paste the computed
# value (not the formula) into the target cell.
$scratch = $ws.Range("Z100")
$scratch.Formula = "=""Oct-24-2023"""
$scratch.Copy()
$ws.Range("B16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Clear()

$ws.Range("D16").Value = "all_submitted_tracker_TM_Oct-24-2023.csv"

$ws.Range("E16").Select()
